{"js": "const body = context.document.body;\n\n// Find the paragraph ending in \"...que har\u00e9.\" (the \"No creo que...\"\n// paragraph) so the new content is anchored robustly rather than via a\n// hard-coded paragraph index.\nconst results = body.search(\"siguiente que har\u00e9.\", { matchCase: false });\nresults.load(\"items\");\nawait context.sync();\n\nconst anchorPara = results.items[0].paragraphs.getFirst();\n\n// Insert a new blank paragraph right after it.\nconst blankPara = anchorPara.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\n// Insert the new sentence about the algorithms used by the author's\n// favorite games/apps (e.g. Google Maps) right after that blank\n// paragraph. The blank paragraph that already existed right before\n// \"No recuerdo haber hecho ...\" is left untouched.\nconst newText =\n  \"Los algoritmos que utilizan mis juegos y aplicaciones favoritas no lo \" +\n  \"s\u00e9 muy bien, pero por ejemplo, cuando uso Google Maps supongo que \" +\n  \"utilizara un algoritmo para encontrar las rutas m\u00e1s r\u00e1pidas.\";\nblankPara.insertParagraph(newText, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Add a new paragraph about algorithms used by games/apps (e.g. Google Maps)\n# right after the \"No creo que ... que har\u00e9.\" paragraph, keeping the\n# pre-existing blank paragraph (before \"No recuerdo haber hecho ...\") intact.\n\n$d = $word.ActiveDocument\n\n# Locate the \"No creo que ... que har\u00e9.\" paragraph (2nd paragraph of the doc)\n# and the blank paragraph that already precedes \"No recuerdo haber hecho ...\".\n$anchorPara = $d.Paragraphs(2)\n$existingBlank = $anchorPara.Next()\n\n# Insert two new (still empty) paragraphs right before that existing blank\n# paragraph: one stays blank, the other will receive the new sentence.\n$existingBlank.Range.InsertParagraphBefore()\n$existingBlank.Range.InsertParagraphBefore()\n\n$newContentPara = $anchorPara.Next().Next()\n\n$ooxml = @'\n<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t xml:space=\"preserve\">Los algoritmos </w:t></w:r><w:r><w:t xml:space=\"preserve\">que utilizan mis juegos y aplicaciones </w:t></w:r><w:r><w:t>favoritas no lo s</w:t></w:r><w:r><w:t>\u00e9</w:t></w:r><w:r><w:t xml:space=\"preserve\"> muy bien, </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>pero</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> por ejemplo</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> cuando uso Google </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Maps</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">supongo </w:t></w:r><w:r><w:t xml:space=\"preserve\">que </w:t></w:r><w:r><w:t>utilizara un algoritmo para encontrar las rutas m\u00e1s</w:t></w:r><w:r><w:t xml:space=\"preserve\"> r\u00e1pidas</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>\n'@\n\n$newContentPara.Range.InsertXML($ooxml)\n"}
